$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
Write-Host "Shape count: " $s.Shapes.Count
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    Write-Host $i $sh.Name $sh.Type
}
